# Daily_Order_Report update — fills in today's received order quantities /
# invoice values per segment (rows 2-9), updates the dealer-level order
# matrix (rows 13-15) with the actual non-zero cells for the day, and adds
# a new dealer row (row 16) for a second "Corporate Territory" order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Segment summary rows (Qty / Invoice Value) ----
# ET
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 28500
# Forklift
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
# IPS
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 22674
# Rickshaw
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 46100
# Solar
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 44102
# Solar P.P
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
# Water
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
# Total
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 141376

# ---- Row 13: Corporate Territory dealer #1 (new dealer for the day) ----
$ws.Range("C13").Value = "Confidence Electric Limited (Factory)"
$ws.Range("D13").Value = "Mojompur, Taltola, Modonpur, Narayangonj"
$ws.Range("E13").Value = "0"
$ws.Range("L13").Value = 0
$ws.Range("AS13").Value = 0
$ws.Range("BB13").Value = 0
$ws.Range("BU13").Value = 2

# ---- Row 14: Gazipur Territory dealer ----
$ws.Range("B14").Value = "Gazipur Territory"
$ws.Range("C14").Value = "M/S Al - Mokka Enterprise (C)"
$ws.Range("D14").Value = "Monipur, Gazipur Sador, Gazipur"
$ws.Range("E14").Value = "01765140095, 01648936899"
$ws.Range("H14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("AC14").Value = 0
$ws.Range("AN14").Value = 2
$ws.Range("BB14").Value = 0

# ---- Row 15: Container Territory dealer ----
$ws.Range("B15").Value = "Container Territory"
$ws.Range("C15").Value = "M/S Bright Renewables Ltd."
$ws.Range("D15").Value = "Sreepur, Mawna, Gazipur"
$ws.Range("E15").Value = ""
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("BI15").Value = 0
$ws.Range("DD15").Value = 2
$ws.Range("DG15").Value = 2
$ws.Range("DH15").Value = 2

# ---- Row 16: new Corporate Territory dealer (added this edit) ----
$ws.Range("A15:DQ15").Copy()
$ws.Range("A16:DQ16").PasteSpecial(-4122)

$ws.Range("A16").Value = 4
$ws.Range("B16").Value = "Corporate Territory"
$ws.Range("C16").Value = "Confidence Group"
$ws.Range("D16").Value = "Unique Trade Centre (UTC), Level 7, 08, Panthapath, Kawran Bazar, Dhaka 1215, Bangladesh"
$ws.Range("E16").Value = ""
$ws.Range("F16:DQ16").Value = 0
$ws.Range("CI16").Value = 2
$ws.Range("CO16").Value = 2
